# Apply "Trade #200 closed" update to live_trading_results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.5       # Current Capital
$summary.Range("B4").Value = -0.49        # Total P&L $
$summary.Range("B6").Value = 200          # Total Trades
$summary.Range("B8").Value = 85           # Losing Trades
$summary.Range("B9").Value = 41.5         # Win Rate %

# ---------------------------------------------------------------------
# Sheet: Strategy Status (volatility_scorer row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.18000000000001   # Capital
$status.Range("D12").Value = 18                  # Trades
$status.Range("E12").Value = -0.82               # P&L $
$status.Range("F12").Value = -0.82               # P&L %
$status.Range("G12").Value = 27.78               # Win Rate %

# ---------------------------------------------------------------------
# Sheet: All Trades - append closed volatility_scorer trade (#200) and
# new open MarketMaking trade (#201)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A201").Value = 200
$allTrades.Range("B201").NumberFormat = "@"
$allTrades.Range("B201").Value = "2026-02-17"
$allTrades.Range("C201").NumberFormat = "@"
$allTrades.Range("C201").Value = "10:08:25"
$allTrades.Range("D201").Value = "volatility_scorer"
$allTrades.Range("E201").Value = "NEUTRAL"
$allTrades.Range("F201").Value = 0.07000000000000001
$allTrades.Range("G201").Value = 0.026357
$allTrades.Range("H201").Value = "CLOSED"
$allTrades.Range("I201").Value = -62.3466
$allTrades.Range("J201").Value = -0.04
$allTrades.Range("K201").Value = 99.18000000000001
$allTrades.Range("L201").Value = 0
$allTrades.Range("M201").Value = 0
$allTrades.Range("N201").Value = 0.85
$allTrades.Range("O201").Value = "Low vol market (score: inf) - ideal for market making"
$allTrades.Range("P201").Value = "early_exit"
$allTrades.Range("Q201").Value = 0.17

$allTrades.Range("A202").Value = 201
$allTrades.Range("B202").NumberFormat = "@"
$allTrades.Range("B202").Value = "2026-02-17"
$allTrades.Range("C202").NumberFormat = "@"
$allTrades.Range("C202").Value = "10:08:25"
$allTrades.Range("D202").Value = "MarketMaking"
$allTrades.Range("E202").Value = "UP"
$allTrades.Range("F202").Value = 0.93
$allTrades.Range("H202").Value = "OPEN"
$allTrades.Range("I202").Value = 0
$allTrades.Range("J202").Value = 0
$allTrades.Range("K202").Value = 100.3271991854616
$allTrades.Range("L202").Value = 0
$allTrades.Range("M202").Value = 0
$allTrades.Range("N202").Value = 0.6
$allTrades.Range("O202").Value = "Normal spread capture: 19600 bps"
$allTrades.Range("Q202").Value = 0

# ---------------------------------------------------------------------
# Sheet: volatility_scorer - append its own closed trade (#200)
# ---------------------------------------------------------------------
$volSheet = $wb.Worksheets.Item("volatility_scorer")

$volSheet.Range("A19").Value = 200
$volSheet.Range("B19").NumberFormat = "@"
$volSheet.Range("B19").Value = "2026-02-17"
$volSheet.Range("C19").NumberFormat = "@"
$volSheet.Range("C19").Value = "10:08:25"
$volSheet.Range("D19").Value = "volatility_scorer"
$volSheet.Range("E19").Value = "NEUTRAL"
$volSheet.Range("F19").Value = 0.07000000000000001
$volSheet.Range("G19").Value = 0.026357
$volSheet.Range("H19").Value = "CLOSED"
$volSheet.Range("I19").Value = -62.3466
$volSheet.Range("J19").Value = -0.04
$volSheet.Range("K19").Value = 99.18000000000001
$volSheet.Range("L19").Value = 0
$volSheet.Range("M19").Value = 0
$volSheet.Range("N19").Value = 0.85
$volSheet.Range("O19").Value = "Low vol market (score: inf) - ideal for market making"
$volSheet.Range("P19").Value = "early_exit"
$volSheet.Range("Q19").Value = 0.17

# ---------------------------------------------------------------------
# Sheet: MarketMaking - append its own new open trade (#201)
# ---------------------------------------------------------------------
$mmSheet = $wb.Worksheets.Item("MarketMaking")

$mmSheet.Range("A184").Value = 201
$mmSheet.Range("B184").NumberFormat = "@"
$mmSheet.Range("B184").Value = "2026-02-17"
$mmSheet.Range("C184").NumberFormat = "@"
$mmSheet.Range("C184").Value = "10:08:25"
$mmSheet.Range("D184").Value = "MarketMaking"
$mmSheet.Range("E184").Value = "UP"
$mmSheet.Range("F184").Value = 0.93
$mmSheet.Range("H184").Value = "OPEN"
$mmSheet.Range("I184").Value = 0
$mmSheet.Range("J184").Value = 0
$mmSheet.Range("K184").Value = 100.3271991854616
$mmSheet.Range("L184").Value = 0
$mmSheet.Range("M184").Value = 0
$mmSheet.Range("N184").Value = 0.6
$mmSheet.Range("O184").Value = "Normal spread capture: 19600 bps"
$mmSheet.Range("Q184").Value = 0
